{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) The \"NE HASZN\u00c1LJ TABLE tag-et!\" bullet moves from the outer list level\n//    (level 0) to the inner list level (level 1) - the same level used by\n//    the other bullets in that sub-list (numId is unchanged).\nlet tablePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"NE HASZN\u00c1LJ TABLE tag-et!\") !== -1) {\n    tablePara = paragraphs.items[i];\n    break;\n  }\n}\nif (!tablePara) {\n  throw new Error('Paragraph containing \"NE HASZN\u00c1LJ TABLE tag-et!\" not found.');\n}\ntablePara.listItem.level = 1;\n\n// 2) Colour the six paragraphs describing the cast/actor rendering rules\n//    green (RGB 008000), matching the rest of the instructional text -\n//    from \"A szerepl\u0151k adatait...\" through \"A sz\u00edn\u00e9szek k\u00e9pei k\u00f6r alak\u00faan...\".\nlet collecting = false;\nlet sawEnd = false;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"A szerepl\u0151k adatait az al\u00e1bbi form\u00e1ban jelen\u00edtsd meg\") !== -1) {\n    collecting = true;\n  }\n  if (collecting) {\n    paragraphs.items[i].font.color = \"#008000\";\n  }\n  if (t.indexOf(\"A sz\u00edn\u00e9szek k\u00e9pei k\u00f6r alak\u00faan\") !== -1) {\n    sawEnd = true;\n    break;\n  }\n}\nif (!sawEnd) {\n  throw new Error(\"Cast-section paragraph range not fully found.\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The \"NE HASZNALJ TABLE tag-et!\" bullet moves from the outer list level\n#    (ilvl 0) to the inner list level (ilvl 1) - same level used by the\n#    other bullets in that sub-list (numId stays the same). COM's\n#    ListLevelNumber is 1-based, so ilvl 0 -> ListLevelNumber 1 and\n#    ilvl 1 -> ListLevelNumber 2.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*NE HASZN\u00c1LJ TABLE tag-et!*\") {\n        $p.Range.ListFormat.ListLevelNumber = 2\n        break\n    }\n}\n\n# 2) Colour the six paragraphs describing the cast/actor rendering rules\n#    green (RGB 008000), matching the rest of the instructional text.\n#    COM colour values are packed as 0x00BBGGRR, so RGB(0,128,0) -> 0x008000.\n$green = 0x008000\n$collecting = $false\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*A szerepl\u0151k adatait az al\u00e1bbi form\u00e1ban jelen\u00edtsd meg*\") {\n        $collecting = $true\n    }\n    if ($collecting) {\n        $p.Range.Font.Color = $green\n    }\n    if ($t -like \"*A sz\u00edn\u00e9szek k\u00e9pei k\u00f6r alak\u00faan*\") {\n        break\n    }\n}\n"}
